# Add a "Name" column between "Date" (D) and "Area" (E), populate it with
# the poster's name for each listing row, and update the Date value for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; this shifts Area/Description/Post URL one column right
# and inherits formatting from the neighboring column, so the header keeps its style.
$ws.Columns("E").Insert()

$ws.Range("E1").Value = "Name"

$names = @{
    2  = "Amyella Hindi"
    3  = "Jianyu Daniel Wu"
    4  = "Amyella Hindi"
    5  = "Bryan Kam"
    6  = "Daniel Cichocki"
    7  = "Laura Scarinci"
    8  = "Eugene Kwon"
    9  = "Jake White"
    10 = "Laura Scarinci"
    11 = "Eun Young Cindy Choi"
    12 = "Sher Khan"
    13 = "Cristopher Rodriguez Ñeke"
    14 = "Lucy Pereira"
    15 = "Gabriel Mairson"
    16 = "Laurene Petitjean"
    17 = "Jack Man"
    18 = "Nikolay Gresko"
    19 = "Jack Man"
    20 = "Peace Lily"
    21 = "Amyella Hindi"
    22 = "Amyella Hindi"
    23 = "Allen Chen"
}

foreach ($row in $names.Keys) {
    $ws.Cells.Item($row, 5).Value = $names[$row]
}

# Update the Date value for row 2
$ws.Range("D2").Value = "1 hr"
